# RF015 - Gerenciar Planos de Capacitacao de TI: wording fixes (1.3 -> 1.4)
#
# The source text used "das Planos ... cadastradas/excluida" (feminine
# agreement, because "Planos" was mistakenly combined with "as") when it
# should agree with "Plano(s)" (masculine). This pass fixes the
# gender/number agreement in the repeated listing/selection/field-fill
# phrases, fixes "um Planos" -> "um Plano" on the selection step, and
# corrects the "nao confirma" listing result so it now correctly states
# the item was NOT removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of exact old text -> new text. Applied as a literal, whole-cell
# replacement everywhere the old text occurs in the used range, so every
# repeated occurrence of each template phrase (one per test case block)
# gets updated consistently.
$replacements = @{
    "SYSTEM exibe a listagem das Planos de Capacitacao de TI cadastradas apenas para visualizacao com a opcao 'Ajuda'" = "SYSTEM exibe a listagem dos Planos de Capacitacao de TI cadastrados apenas para visualizacao com a opcao 'Ajuda'"
    "SYSTEM exibe a listagem das Planos de Capacitacao de TI cadastradas com opcoes de 'Novo', 'Editar', 'Excluir' e 'Ajuda'" = "SYSTEM exibe a listagem dos Planos de Capacitacao de TI cadastrados com opcoes de 'Novo', 'Editar', 'Excluir' e 'Ajuda'"
    "Lider de Pessoas seleciona um Planos de Capacitacao de TI da listagem" = "Lider de Pessoas seleciona um Plano de Capacitacao de TI da listagem"
    "SYSTEM exibe a listagem das Planos de Capacitacao de TI com a Capacitacao de TI excluida" = "SYSTEM exibe a listagem dos Planos de Capacitacao de TI com a Capacitacao de TI nao excluida"
    "SYSTEM exibe a listagem das Planos de Capacitacao de TI sem a Capacitacao de TI excluida" = "SYSTEM exibe a listagem dos Planos de Capacitacao de TI sem a Capacitacao de TI excluida"
    "Lider de Pessoas escolha o 'Periodo Avaliativo' apropriado no campo de selecao" = "Lider de Pessoas escolhe o 'Periodo Avaliativo' apropriado no campo de selecao"
    "Lider de Pessoas selecione a 'Unidade' correspondente no campo de selecao de unidade" = "Lider de Pessoas seleciona a 'Unidade' correspondente no campo de selecao de unidade"
    "Lider de Pessoas preencha o campo 'Possiveis Capacitacoes' com informacoes sobre capacitacoes adicionais" = "Lider de Pessoas preenche o campo 'Possiveis Capacitacoes' com informacoes sobre capacitacoes adicionais"
    "Lider de Pessoas preencha o campo 'Observacao' com informacoes adicionais ou relevantes sobre o plano de capacitacao" = "Lider de Pessoas preenche o campo 'Observacao' com informacoes adicionais ou relevantes sobre o plano de capacitacao"
}

$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($null -ne $val -and $replacements.ContainsKey($val)) {
            $cell.Value = $replacements[$val]
        }
    }
}
